$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) figures.
# Cells in column D whose new text looks like a plain number need their
# NumberFormat forced to Text ("@") first, otherwise Excel would silently
# convert the assigned string into a floating point number and mangle
# formatting such as trailing zeros (e.g. "71.00" -> 71).

$ws.Range("D2").Value = "38.630.13"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").Value = "2.093.26"
$ws.Range("E3").Value = "  +2.70%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.65"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.14"
$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +1.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0842"
$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").Value = "2.398.45"
$ws.Range("E12").Value = "  +2.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.85"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.35"
$ws.Range("E14").Value = "  +6.14%  "

$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("E16").Value = "  +4.97%  "

$ws.Range("D17").Value = "2.110.08"
$ws.Range("E17").Value = "  +3.38%  "

$ws.Range("D18").Value = "38.537.35"
$ws.Range("E18").Value = "  +1.96%  "

$ws.Range("E19").Value = "  +3.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.00"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.43"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +1.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.23"
$ws.Range("E26").Value = "  +1.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.45"
$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.136"
$ws.Range("E28").Value = "  +5.52%  "

$ws.Range("E29").Value = "  +1.69%  "

$ws.Range("E30").Value = "  +7.03%  "

$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("E32").Value = "  +4.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.76"
$ws.Range("E33").Value = "  +5.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  +2.32%  "

$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.55"
$ws.Range("E36").Value = "  +1.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  +3.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.50"
$ws.Range("E40").Value = "  +2.65%  "

$ws.Range("D41").Value = "1.544.21"
$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.15"
$ws.Range("E42").Value = "  +4.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0220"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("E44").Value = "  +1.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0921"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.74"
$ws.Range("E46").Value = "  +10.03%  "

$ws.Range("E47").Value = "  +2.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.03"
$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  +0.63%  "

$ws.Range("D51").Value = "2.290.15"
$ws.Range("E51").Value = "  +2.78%  "
